$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Score corrections across several students' rows
$ws.Range("B2").Value = 4
$ws.Range("C3").Value = 4
$ws.Range("C7").Value = 4
$ws.Range("B9").Value = 4
$ws.Range("D10").Value = 5
$ws.Range("B11").Value = 4
$ws.Range("E16").Value = 5
$ws.Range("B19").Value = 5
$ws.Range("D19").Value = 0
$ws.Range("D20").Value = 5
$ws.Range("B21").Value = 5
$ws.Range("D23").Value = 0
$ws.Range("D24").Value = 5
$ws.Range("C25").Value = 4

# Update notes for students whose work is now fully correct
$ws.Range("L3").Value = "переписана верно все номера"
$ws.Range("L7").Value = "переписана верно все номера"
$ws.Range("L25").Value = "переписана верно все номера"

# Restore active cell selection
$ws.Range("B10").Select()
